$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 1
    6  = 0
    7  = 4
    8  = 2
    9  = 2
    10 = 0
    11 = 3
    12 = 0
    13 = 3
    14 = 0
    15 = 0
    16 = 1
    17 = 1
    18 = 2
    19 = 0
    20 = 1
    21 = 1
    22 = 2
    23 = 2
    24 = 1
    25 = 1
    26 = 2
    27 = 2
    28 = 1
    29 = 1
    30 = 1
    31 = 1
    32 = 1
    33 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
